$wb = $excel.ActiveWorkbook
$wpWs = $wb.Worksheets.Item("Weekly Points")
$cpWs = $wb.Worksheets.Item("Cumulative Points")

# ---- Weekly Points: Gameweek 6 rows (152-181) ----
$data = @(
    @(6,"Bazzers Ballers","Onana","Manchester United","Goalkeeper",3),
    @(6,"Bazzers Ballers","Areola","West Ham","Goalkeeper",0),
    @(6,"Bazzers Ballers","Van de Ven","Tottenham","Defender",11),
    @(6,"Bazzers Ballers","Andersen","Fulham","Defender",8),
    @(6,"Bazzers Ballers","White","Arsenal","Defender",0),
    @(6,"Bazzers Ballers","Romero","Tottenham","Defender",6),
    @(6,"Bazzers Ballers","Johnson","Ipswich","Defender",0),
    @(6,"Bazzers Ballers","Luis Diaz","Liverpool","Midfield",2),
    @(6,"Bazzers Ballers","Thomas","Arsenal","Midfield",2),
    @(6,"Bazzers Ballers","Bowen","West Ham","Midfield",2),
    @(6,"Bazzers Ballers","Palmer","Chelsea","Midfield",25),
    @(6,"Bazzers Ballers","Onana","Aston Villa","Midfield",2),
    @(6,"Bazzers Ballers","Haaland","Manchester City","Forward",4),
    @(6,"Bazzers Ballers","Joao Pedro","Brighton","Forward",0),
    @(6,"Bazzers Ballers","Mateta","Crystal Palace","Forward",2),
    @(6,"WHU-Tang-Clan","Arrizaballago","Bournemouth","Goalkeeper",2),
    @(6,"WHU-Tang-Clan","Areola","West Ham","Goalkeeper",0),
    @(6,"WHU-Tang-Clan","White","Arsenal","Defender",0),
    @(6,"WHU-Tang-Clan","Cucurello","Chelsea","Defender",0),
    @(6,"WHU-Tang-Clan","Romero","Tottenham","Defender",6),
    @(6,"WHU-Tang-Clan","Murillo","Nottingham Forest","Defender",0),
    @(6,"WHU-Tang-Clan","De Ligt","Manchester United","Defender",1),
    @(6,"WHU-Tang-Clan","Bowen","West Ham","Midfield",2),
    @(6,"WHU-Tang-Clan","Palmer","Chelsea","Midfield",25),
    @(6,"WHU-Tang-Clan","Amad","Manchester United","Midfield",1),
    @(6,"WHU-Tang-Clan","Smith-Rowe","Fulham","Midfield",3),
    @(6,"WHU-Tang-Clan","Madueke","Chelsea","Midfield",0),
    @(6,"WHU-Tang-Clan","Wellbeck","Brighton","Forward",2),
    @(6,"WHU-Tang-Clan","Haaland","Manchester City","Forward",4),
    @(6,"WHU-Tang-Clan","Havertz","Arsenal","Forward",6)
)

$startRow = 152
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $wpWs.Cells.Item($r, 1).Value = $row[0]
    $wpWs.Cells.Item($r, 2).Value = $row[1]
    $wpWs.Cells.Item($r, 3).Value = $row[2]
    $wpWs.Cells.Item($r, 4).Value = $row[3]
    $wpWs.Cells.Item($r, 5).Value = $row[4]
    $wpWs.Cells.Item($r, 6).Value = $row[5]
}

# ---- Cumulative Points: Gameweek 6 rows (14-15) ----
$cpWs.Cells.Item(14, 1).Value = "Bazzers Ballers"
$cpWs.Cells.Item(14, 2).Value = 6
$cpWs.Range("C14").Formula = "=SUM('Weekly Points'!F152:F166)+C12"

$cpWs.Cells.Item(15, 1).Value = "WHU-Tang-Clan"
$cpWs.Cells.Item(15, 2).Value = 6
$cpWs.Range("C15").Formula = "=SUM('Weekly Points'!F167:F181)+C13"

$excel.Calculate()

# ---- View state updates ----
$cpWs.Activate()
$cpWs.Range("B16").Select()

$wpWs.Activate()
$wpWs.Range("H153").Select()
